# Renames the placeholder interface name "YYY"/"yyy" to "interfaceRubrica"
# throughout the document. Each occurrence is textually embedded inside a
# larger run, so after changing the text we "stamp" a run-boundary around
# the replacement by toggling Bold on/off (a no-op formatting-wise) which
# forces Word to split the paragraph into separate <w:r> elements at the
# boundary, matching how Word itself splits a run when you replace only
# part of its text via Find & Replace / manual selection.

function Split-RunBoundary($range) {
    $range.Bold = 1
    $range.Bold = 0
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Target 1: "……YYY" -> "……" + "interfaceRubrica"   (color FF0000 run)
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("……YYY", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Target 1 text not found" }
$matchStart = $r.Start

$oldStart = $matchStart + 2          # skip the leading "……"
$oldEnd = $matchStart + 5            # "YYY" is 3 chars

$rOld = $d.Range($oldStart, $oldEnd)
$rOld.Text = "interfaceRubrica"

$newLen = "interfaceRubrica".Length
$rNew = $d.Range($oldStart, $oldStart + $newLen)
Split-RunBoundary $rNew

# ---------------------------------------------------------------------
# Target 2: "A class that implements the interface yyy must provide an
# implementation for Contatto cerca(String cognome). "
# -> "...interface " + "interfaceRubrica" + " must provide..."
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("A class that implements the interface yyy must provide an implementation for Contatto cerca(String cognome).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Target 2 text not found" }
$matchStart = $r.Start

$prefix = "A class that implements the interface "
$oldWord = "yyy"
$oldStart = $matchStart + $prefix.Length
$oldEnd = $oldStart + $oldWord.Length

$rOld = $d.Range($oldStart, $oldEnd)
$rOld.Text = "interfaceRubrica"

$newLen = "interfaceRubrica".Length
$rNew = $d.Range($oldStart, $oldStart + $newLen)
Split-RunBoundary $rNew

# ---------------------------------------------------------------------
# Target 3: "Any concrete class that implements the YYY interface must
# define a Contatto cerca(String cognome) instance method. "
# -> "...implements the " + "interfaceRubrica" + " interface must..."
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Any concrete class that implements the YYY interface must define a Contatto cerca(String cognome) instance method.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Target 3 text not found" }
$matchStart = $r.Start

$prefix = "Any concrete class that implements the "
$oldWord = "YYY"
$oldStart = $matchStart + $prefix.Length
$oldEnd = $oldStart + $oldWord.Length

$rOld = $d.Range($oldStart, $oldEnd)
$rOld.Text = "interfaceRubrica"

$newLen = "interfaceRubrica".Length
$rNew = $d.Range($oldStart, $oldStart + $newLen)
Split-RunBoundary $rNew

# ---------------------------------------------------------------------
# Target 4: "YYY{" -> "interfaceRubrica" + "{"
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("YYY{", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Target 4 text not found" }
$matchStart = $r.Start

$oldWord = "YYY"
$oldStart = $matchStart
$oldEnd = $oldStart + $oldWord.Length

$rOld = $d.Range($oldStart, $oldEnd)
$rOld.Text = "interfaceRubrica"

$newLen = "interfaceRubrica".Length
$rNew = $d.Range($oldStart, $oldStart + $newLen)
Split-RunBoundary $rNew

Write-Output "Replacements complete."
